# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1404
$ws1.Range("F6").Value  = 70
$ws1.Range("F7").Value  = 11891
$ws1.Range("F8").Value  = 4438
$ws1.Range("F13").Value = 2567
$ws1.Range("F14").Value = 1107
$ws1.Range("F15").Value = 164
$ws1.Range("F17").Value = 5158
$ws1.Range("F19").Value = 194
$ws1.Range("F20").Value = 534
$ws1.Range("F21").Value = 11386
$ws1.Range("F22").Value = 11372
$ws1.Range("F24").Value = 51
$ws1.Range("F27").Value = 51

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1404
$ws4.Range("F6").Value  = 70
$ws4.Range("F7").Value  = 11891
$ws4.Range("F8").Value  = 4438
$ws4.Range("F13").Value = 2567
$ws4.Range("F15").Value = 1107
$ws4.Range("F16").Value = 164
$ws4.Range("F18").Value = 5158
$ws4.Range("F20").Value = 194
$ws4.Range("F21").Value = 534
$ws4.Range("F22").Value = 11386
$ws4.Range("F23").Value = 11372
$ws4.Range("F25").Value = 51
$ws4.Range("F28").Value = 51

$wb.Save()
